$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Rows.Insert() in this runtime shifts cell values/formulas correctly
# but does NOT shift the worksheet's <hyperlink> anchor refs for rows below
# the insertion point. Work around it by wiping all hyperlinks up-front and
# re-creating every one of them (old + new) at their correct final address
# after all row inserts are done.
$ws.Range("A1").Hyperlinks.Delete() | Out-Null

# --- Insert two new rows before the first table's Subtotal row (row 7) ---
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# --- Insert two new rows before the second table's Subtotal row ---
# After the first pair of inserts above, the second table's subtotal row
# (originally row 15) is now row 17.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

function Fill-NutScrewRows($hexRow, $screwRow) {
    # Row with the hex nut (wraps to two lines, so it is taller)
    $ws.Rows.Item($hexRow).RowHeight = 28.8

    $ws.Cells.Item($hexRow, 1).Value = 6
    $ws.Cells.Item($hexRow, 2).Value = 4
    $ws.Cells.Item($hexRow, 3).Value = "36-4701-ND"
    $ws.Cells.Item($hexRow, 4).Value = 4701
    $ws.Cells.Item($hexRow, 5).Value = "`t`nHEX NUT 0.245`" STEEL 6-32"
    $ws.Cells.Item($hexRow, 5).WrapText = $true
    $ws.Cells.Item($hexRow, 6).Value = 0.1
    $ws.Cells.Item($hexRow, 8).Value = "https://www.digikey.com/en/products/detail/keystone-electronics/4701/316272"

    # Row with the machine screw
    $ws.Cells.Item($screwRow, 1).Value = 7
    $ws.Cells.Item($screwRow, 2).Value = 4
    $ws.Cells.Item($screwRow, 3).Value = "36-9309-ND"
    $ws.Cells.Item($screwRow, 4).Value = 9309
    $ws.Cells.Item($screwRow, 5).Value = "MACH SCREW PAN HEAD SLOTTED 6-32"
    $ws.Cells.Item($screwRow, 6).Value = 0.1
    $ws.Cells.Item($screwRow, 8).Value = "https://www.digikey.com/en/products/detail/keystone-electronics/9309/2746088"
}

Fill-NutScrewRows 7 8
Fill-NutScrewRows 17 18

# --- Fix up the subtotal formulas to include the two new rows ---
$ws.Range("G9").Formula = "=(F2*B2+F3*B3+F4*B4+F5*B5+F6*B6+F7*B7+F8*B8)"
$ws.Range("G19").Formula = "=(F12*B12+F13*B13+F14*B14+F15*B15+F16*B16+F17*B17+F18*B18)"

# --- Re-create every hyperlink (both pre-existing ones that shifted down
#     and the four brand-new ones) now that all rows are in their final
#     positions. ---
$ws.Hyperlinks.Add($ws.Range("H2"),  "https://www.digikey.com/product-detail/en/adafruit-industries-llc/2995/1528-1562-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H3"),  "https://www.digikey.com/products/en?keywords=COM-12986") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H4"),  "https://www.digikey.com/product-detail/en/adafruit-industries-llc/2886/1528-1560-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H5"),  "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/CF14JT470R/CF14JT470RCT-ND/1830342") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H6"),  "https://www.digikey.com/product-detail/en/cui-inc/SJ1-3535NG/CP1-3535NG-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H7"),  "https://www.digikey.com/en/products/detail/keystone-electronics/4701/316272") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H8"),  "https://www.digikey.com/en/products/detail/keystone-electronics/9309/2746088") | Out-Null

$ws.Hyperlinks.Add($ws.Range("H12"), "https://www.digikey.com/product-detail/en/adafruit-industries-llc/2829/1528-1517-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H13"), "https://www.digikey.com/products/en?keywords=COM-12986") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H14"), "https://www.digikey.com/product-detail/en/adafruit-industries-llc/2886/1528-1560-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H15"), "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/CF14JT470R/CF14JT470RCT-ND/1830342") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H16"), "https://www.digikey.com/product-detail/en/cui-inc/SJ1-3535NG/CP1-3535NG-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H17"), "https://www.digikey.com/en/products/detail/keystone-electronics/4701/316272") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H18"), "https://www.digikey.com/en/products/detail/keystone-electronics/9309/2746088") | Out-Null

# --- Update selection to match the final state ---
$ws.Range("G19").Select()

Write-Output "done"
